$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "New device with an ADC error.  Needs ADCON2 = 15 ' Set Negative Reference Setting to ADNREF in ADCON1 `nSee https://sourceforge.net/p/gcbasic/discussion/629990/thread/9b69d693/#e018"

$ws.Range("A16").Value = 15
$ws.Range("A16").HorizontalAlignment = -4131
$ws.Range("A16").VerticalAlignment = -4160
$ws.Range("A16").WrapText = $false

$ws.Range("B16").Value = "OPEN"
$ws.Range("B16").HorizontalAlignment = -4131
$ws.Range("B16").VerticalAlignment = -4160
$ws.Range("B16").WrapText = $false

$ws.Range("D16").Value = $newText
$ws.Range("D16").HorizontalAlignment = -4131
$ws.Range("D16").VerticalAlignment = -4160
$ws.Range("D16").WrapText = $true

$ws.Rows.Item(16).RowHeight = 45

$ws.Range("D1").Select()

$wb.Save()
